# Atualização automática via cronjob
# Remove the old "2025-04-02" rows (original rows 2-8), shifting remaining
# rows up, then refresh the "quantidade_atipica" (A) and "estoque_atualizado" (G)
# values for the rows that remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 7 rows that belonged to the 2025-04-02 batch (rows 2 through 8)
$ws.Range("A2:H8").EntireRow.Delete()

# Updated quantidade_atipica (A) / estoque_atualizado (G) values for the
# rows that remain after the shift (now rows 2-9)
$updates = @(
    @{ Row = 2; A = 5;  G = 174 },
    @{ Row = 3; A = 1;  G = 29 },
    @{ Row = 4; A = 2;  G = -20 },
    @{ Row = 5; A = 3;  G = 493 },
    @{ Row = 6; A = 4;  G = 207 },
    @{ Row = 7; A = 6;  G = 86 },
    @{ Row = 8; A = 7;  G = 95 },
    @{ Row = 9; A = 0;  G = 134 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}
